# Scheduled-runner update to Sheets/Sephirot_Profits.xlsx (workbook with
# per-job ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets). Refreshes the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for the
# rows whose underlying Universalis market data moved since the last run.
# Only numeric literal values change; no formulas, formatting, or
# structural edits are involved.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 17000
$ws.Range("J38").Value = 17000
$ws.Range("L38").Value = 51000
$ws.Range("N38").Value = -51744
$ws.Range("H100").Value = 1529.2858
$ws.Range("I100").Value = 1452
$ws.Range("J100").Value = 1632.3334
$ws.Range("K100").Value = 1452
$ws.Range("L100").Value = 1632.3334
$ws.Range("M100").Value = -911
$ws.Range("N100").Value = -2714.3334
$ws.Range("H112").Value = 3566.2
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3566.2
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 10698.6
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = -12914.6
$ws.Range("H118").Value = 464.14285
$ws.Range("I118").Value = 464.14285
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1392.42855
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 264.5714499999999
$ws.Range("N118").Value = $null
$ws.Range("H132").Value = 1389.9474
$ws.Range("I132").Value = 1244.9445
$ws.Range("K132").Value = 3734.8335
$ws.Range("M132").Value = -1204.8335
$ws.Range("H138").Value = 4222.6816
$ws.Range("J138").Value = 4178.8945
$ws.Range("L138").Value = 12536.6835
$ws.Range("N138").Value = -22816.6835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5756.9546
$ws.Range("I32").Value = 5078.7144
$ws.Range("K32").Value = 5078.7144
$ws.Range("M32").Value = -4791.7144
$ws.Range("H45").Value = 3250
$ws.Range("I45").Value = 2500
$ws.Range("K45").Value = 2500
$ws.Range("M45").Value = -2123
$ws.Range("H61").Value = 2965
$ws.Range("I61").Value = 1439
$ws.Range("K61").Value = 1439
$ws.Range("M61").Value = -1227
$ws.Range("H74").Value = 935.5
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 935.5
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null
$ws.Range("H132").Value = 1831.8846
$ws.Range("I132").Value = 881.55
$ws.Range("K132").Value = 2644.65
$ws.Range("M132").Value = -114.6499999999996
$ws.Range("H136").Value = 2965
$ws.Range("I136").Value = 1439
$ws.Range("K136").Value = 4317
$ws.Range("M136").Value = -1767

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2119.75
$ws.Range("I80").Value = 80.666664
$ws.Range("J80").Value = 3343.2
$ws.Range("K80").Value = 80.666664
$ws.Range("L80").Value = 3343.2
$ws.Range("M80").Value = 917.333336
$ws.Range("N80").Value = -5339.2
$ws.Range("H83").Value = 2119.75
$ws.Range("I83").Value = 80.666664
$ws.Range("J83").Value = 3343.2
$ws.Range("K83").Value = 403.33332
$ws.Range("L83").Value = 16716
$ws.Range("M83").Value = 4588.66668
$ws.Range("N83").Value = -26700
$ws.Range("H86").Value = 2740.875
$ws.Range("I86").Value = 1866
$ws.Range("J86").Value = 4199
$ws.Range("K86").Value = 1866
$ws.Range("L86").Value = 4199
$ws.Range("M86").Value = -743
$ws.Range("N86").Value = -6445
$ws.Range("H89").Value = 2740.875
$ws.Range("I89").Value = 1866
$ws.Range("J89").Value = 4199
$ws.Range("K89").Value = 9330
$ws.Range("L89").Value = 20995
$ws.Range("M89").Value = -3714
$ws.Range("N89").Value = -32227
$ws.Range("H134").Value = 8294.947
$ws.Range("I134").Value = 2980.5
$ws.Range("K134").Value = 8941.5
$ws.Range("M134").Value = -6406.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3668.0938
$ws.Range("I31").Value = 2818.0952
$ws.Range("J31").Value = 5290.8184
$ws.Range("K31").Value = 2818.0952
$ws.Range("L31").Value = 5290.8184
$ws.Range("M31").Value = -2523.0952
$ws.Range("N31").Value = -5880.8184
$ws.Range("H34").Value = 3668.0938
$ws.Range("I34").Value = 2818.0952
$ws.Range("J34").Value = 5290.8184
$ws.Range("K34").Value = 2818.0952
$ws.Range("L34").Value = 5290.8184
$ws.Range("M34").Value = -2616.0952
$ws.Range("N34").Value = -5694.8184
$ws.Range("H62").Value = 4288.5713
$ws.Range("I62").Value = 3753.75
$ws.Range("J62").Value = 5001.6665
$ws.Range("K62").Value = 3753.75
$ws.Range("L62").Value = 5001.6665
$ws.Range("M62").Value = -3129.75
$ws.Range("N62").Value = -6249.6665
$ws.Range("H65").Value = 4288.5713
$ws.Range("I65").Value = 3753.75
$ws.Range("J65").Value = 5001.6665
$ws.Range("K65").Value = 18768.75
$ws.Range("L65").Value = 25008.3325
$ws.Range("M65").Value = -15648.75
$ws.Range("N65").Value = -31248.3325
$ws.Range("H70").Value = 100108
$ws.Range("J70").Value = 100108
$ws.Range("L70").Value = 100108
$ws.Range("N70").Value = -100738
$ws.Range("H73").Value = 100108
$ws.Range("J73").Value = 100108
$ws.Range("L73").Value = 100108
$ws.Range("N73").Value = -102292
$ws.Range("H134").Value = 4513.143
$ws.Range("I134").Value = 4432.5
$ws.Range("J134").Value = 4997
$ws.Range("K134").Value = 13297.5
$ws.Range("L134").Value = 14991
$ws.Range("M134").Value = -10762.5
$ws.Range("N134").Value = -20061
$ws.Range("H141").Value = 63328.668
$ws.Range("J141").Value = 63328.668
$ws.Range("L141").Value = 63328.668
$ws.Range("N141").Value = -73688.66800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 9664
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 9664
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 28992
$ws.Range("M100").Value = $null
$ws.Range("N100").Value = -30614
$ws.Range("H113").Value = 1254.8
$ws.Range("I113").Value = 2669
$ws.Range("K113").Value = 8007
$ws.Range("M113").Value = -5837
$ws.Range("H131").Value = 2335.5715
$ws.Range("I131").Value = 2430
$ws.Range("K131").Value = 7290
$ws.Range("M131").Value = -2250
$ws.Range("H141").Value = 995
$ws.Range("I141").Value = 995
$ws.Range("K141").Value = 2985
$ws.Range("M141").Value = 2195

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2366.3333
$ws.Range("I102").Value = 2366.3333
$ws.Range("K102").Value = 2366.3333
$ws.Range("M102").Value = -744.3332999999998
$ws.Range("H113").Value = 2546.6
$ws.Range("I113").Value = 2433.25
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2433.25
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -263.25
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1674.75
$ws.Range("I68").Value = 1674.75
$ws.Range("K68").Value = 1674.75
$ws.Range("M68").Value = -925.75
$ws.Range("H71").Value = 1674.75
$ws.Range("I71").Value = 1674.75
$ws.Range("K71").Value = 8373.75
$ws.Range("M71").Value = -4629.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3290.795
$ws.Range("I132").Value = 1975.1154
$ws.Range("K132").Value = 5925.3462
$ws.Range("M132").Value = -3395.3462
$ws.Range("H136").Value = 2800.5
$ws.Range("I136").Value = 2575.75
$ws.Range("K136").Value = 7727.25
$ws.Range("M136").Value = -5177.25

